$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A91").Value = "16-12-2025"
$ws.Range("B91").Value = "The price of gold in India today is ₹13,386 per gram for 24 karat gold, ₹12,270 per gram for 22 karat gold and ₹10,039 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A90").Copy()
$ws.Range("A91").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B90").Copy()
$ws.Range("B91").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
